$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 'valid_min' / 'valid_max' attribute rows from each of the three
# variable blocks (height_above_snow_surface, wind_speed, wind_from_direction).
# Deleting from the bottom block upward keeps the earlier row numbers valid.
$ws.Rows("33:34").Delete()
$ws.Rows("21:22").Delete()
$ws.Rows("9:10").Delete()

# Restore the view state recorded after the edit.
$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("A29:C30").Select()
$ws.Range("C30").Activate()
